$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) column cells are treated as text so values like "560.07"
# are not auto-converted to numbers (matching original inlineStr text cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.863.90"
$ws.Range("E2").Value = "  -0.78%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.398.35"
$ws.Range("E3").Value = "  -1.13%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.07"
$ws.Range("E5").Value = "  +0.57%  "

$ws.Range("E6").Value = "  -1.17%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  -0.70%  "

$ws.Range("E9").Value = "  -1.94%  "

$ws.Range("E10").Value = "  -1.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.24"
$ws.Range("E11").Value = "  -2.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.45"
$ws.Range("E13").Value = "  -3.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000171"
$ws.Range("E14").Value = "  -2.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.831.34"
$ws.Range("E15").Value = "  -1.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.744.79"
$ws.Range("E16").Value = "  -0.59%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.387.45"
$ws.Range("E17").Value = "  -1.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.17"
$ws.Range("E18").Value = "  +0.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "320.42"
$ws.Range("E19").Value = "  -1.54%  "

$ws.Range("E20").Value = "  -1.25%  "

$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.52"
$ws.Range("E23").Value = "  +0.69%  "

$ws.Range("E24").Value = "  -3.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.76"
$ws.Range("E25").Value = "  -4.54%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "560.68"
$ws.Range("E26").Value = "  -2.68%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.515.40"
$ws.Range("E28").Value = "  -1.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0930"
$ws.Range("E29").Value = "  -2.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.16"
$ws.Range("E30").Value = "  -3.13%  "

$ws.Range("E31").Value = "  -5.20%  "

$ws.Range("E32").Value = "  -1.59%  "

$ws.Range("E33").Value = "  -0.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.49"

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("E36").Value = "  -2.39%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "152.47"
$ws.Range("E37").Value = "  +2.68%  "

$ws.Range("E38").Value = "  -5.98%  "

$ws.Range("E39").Value = "  -2.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.48"
$ws.Range("E40").Value = "  -2.04%  "

$ws.Range("E41").Value = "  -6.02%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("E43").Value = "  -3.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "147.26"
$ws.Range("E44").Value = "  -3.43%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.59"
$ws.Range("E45").Value = "  -1.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0527"
$ws.Range("E46").Value = "  -3.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.74"
$ws.Range("E47").Value = "  -3.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.584"
$ws.Range("E48").Value = "  -1.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0916"
$ws.Range("E49").Value = "  +0.14%  "

$ws.Range("E50").Value = "  -2.15%  "

$ws.Range("E51").Value = "  +0.39%  "
